# Update automatico via Actualizar 07-01-2020 06-20-54
# Append two new daily-report rows (2020-06-29 and 2020-06-30) to the
# "Condicion_Pacientes" table on sheet "Hoja1", then expand the table /
# autofilter to cover them and refresh the active view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Copy the formatting of the last existing data row (108) down into the two
# new rows so the date column keeps its date number format and the numeric
# columns keep their centered alignment.
$ws.Range("A108:F108").Copy()
$ws.Range("A109:F109").PasteSpecial(-4122)
$ws.Range("A108:F108").Copy()
$ws.Range("A110:F110").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 109 -> 29/06/2020 (serial date 44011)
$ws.Cells.Item(109, 1).Value = 44011
$ws.Cells.Item(109, 2).Value = 1540
$ws.Cells.Item(109, 3).Value = 736
$ws.Cells.Item(109, 4).Value = 758
$ws.Cells.Item(109, 5).Value = 399
$ws.Cells.Item(109, 6).Value = 56

# Row 110 -> 30/06/2020 (serial date 44012)
$ws.Cells.Item(110, 1).Value = 44012
$ws.Cells.Item(110, 2).Value = 1697
$ws.Cells.Item(110, 3).Value = 740
$ws.Cells.Item(110, 4).Value = 730
$ws.Cells.Item(110, 5).Value = 415
$ws.Cells.Item(110, 6).Value = 56

# Grow the "Condicion_Pacientes" table (and its autofilter) to include rows
# 109:110, same as resizing ref from A1:F108 to A1:F110.
$table = $ws.ListObjects.Item("Condicion_Pacientes")
$table.Resize($ws.Range("A1:F110"))

# Reflect the newly active area in the sheet view / selection.
$ws.Activate()
$ws.Range("D110").Select()
